$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently holds the text "R40". The rule row is being re-keyed to
# "1" (still a text label, not a number) while everything else about the row
# (styling, the numeric From/To bounds in C11/D11, the Good Night greeting in
# E11, etc.) stays the same. Force the cell to remain text so "1" isn't
# reinterpreted as a numeric value.
$b11 = $ws.Range("B11")
$b11.NumberFormat = "@"
$b11.Value = "1"
